# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Espinaca" at Feria Lagunitas de
# Puerto Montt. The new observation belongs right after the existing
# row 21 (it sorts there chronologically within the sheet's ordering),
# so insert a blank row at row 22 - shifting the former rows 22:37 down
# to 23:38 - and populate the new row with the week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 22:37 down to 23:38, leaving a blank row 22 (inherits the
# formatting, including the date-column style, from the row above).
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record.
$ws.Cells.Item(22, 1).Value  = 4
$ws.Cells.Item(22, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(22, 3).Value  = "Los Lagos"
$ws.Cells.Item(22, 4).Value  = 44789
$ws.Cells.Item(22, 5).Value  = 10
$ws.Cells.Item(22, 6).Value  = 100112012
$ws.Cells.Item(22, 7).Value  = "Espinaca"
$ws.Cells.Item(22, 8).Value  = "Sin especificar"
$ws.Cells.Item(22, 9).Value  = "Primera"
$ws.Cells.Item(22, 10).Value = 30
$ws.Cells.Item(22, 11).Value = 15000
$ws.Cells.Item(22, 12).Value = 15000
$ws.Cells.Item(22, 13).Value = 15000
$ws.Cells.Item(22, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(22, 15).Value = "Región Metropolitana"
$ws.Cells.Item(22, 16).Value = 1500
$ws.Cells.Item(22, 17).Value = 10
$ws.Cells.Item(22, 18).Value = "Hortaliza"
